$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Previously (pre ClosedXML 0.100) a cell could hold a formula-like string
# such as &="Total: "<<sum>> and the template engine would read the literal
# text of the formula. That no longer works, so the template is fixed to
# just contain the tag itself as plain text.
$ws.Range("G6").Value = "<<sum>>"
